$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("B18").Value = 27078
$ws.Range("C18").Value = 346
$ws.Range("E18").Value = 9890
# Row 19
$ws.Range("B19").Value = 19022
$ws.Range("C19").Value = 181
$ws.Range("D19").Value = 519
$ws.Range("E19").Value = 17846
$ws.Range("F19").Value = 222
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 657
# Row 34
$ws.Range("E34").Value = 6715
$ws.Range("G34").Value = 6
$ws.Range("H34").Value = 158
# Row 44
$ws.Range("A44").Value = "Singapur"
$ws.Range("B44").Value = 5050
$ws.Range("C44").Value = 623
$ws.Range("D44").Value = 683
$ws.Range("E44").Value = 4357
$ws.Range("F44").Value = 29
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 10
# Row 45
$ws.Range("A45").Value = "Bielorrusia"
$ws.Range("B45").Value = 4779
$ws.Range("C45").Value = 575
$ws.Range("D45").Value = 342
$ws.Range("E45").Value = 4395
$ws.Range("F45").Value = 65
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 42
# Row 46
$ws.Range("A46").Value = "Catar"
$ws.Range("B46").Value = 4663
$ws.Range("C46").Value = 560
$ws.Range("D46").Value = 464
$ws.Range("E46").Value = 4192
$ws.Range("F46").Value = 37
$ws.Range("H46").Value = 7
# Row 47
$ws.Range("A47").Value = "Ucrania"
$ws.Range("B47").Value = 4662
$ws.Range("C47").Value = 501
$ws.Range("D47").Value = 246
$ws.Range("E47").Value = 4291
$ws.Range("F47").Value = 45
$ws.Range("G47").Value = 9
$ws.Range("H47").Value = 125
# Row 66
$ws.Range("B66").Value = 1658
$ws.Range("C66").Value = 134
$ws.Range("D66").Value = 258
$ws.Range("E66").Value = 1395
$ws.Range("G66").Value = 2
$ws.Range("H66").Value = 5
# Row 67
$ws.Range("D67").Value = 306
$ws.Range("E67").Value = 1157
# Row 75
$ws.Range("B75").Value = 1199
$ws.Range("C75").Value = 32
$ws.Range("D75").Value = 320
$ws.Range("E75").Value = 833
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = 46
# Row 79
$ws.Range("A79").Value = "Eslovaquia"
$ws.Range("B79").Value = 1049
$ws.Range("C79").Value = 72
$ws.Range("D79").Value = 167
$ws.Range("E79").Value = 874
$ws.Range("F79").Value = 5
$ws.Range("H79").Value = 8
# Row 80
$ws.Range("A80").Value = "Hong Kong"
$ws.Range("B80").Value = 1022
$ws.Range("C80").Value = 4
$ws.Range("D80").Value = 533
$ws.Range("E80").Value = 485
$ws.Range("F80").Value = 9
$ws.Range("H80").Value = 4
# Row 81
$ws.Range("A81").Value = "Camerun"
$ws.Range("B81").Value = 996
$ws.Range("D81").Value = 164
$ws.Range("E81").Value = 810
$ws.Range("F81").Value = 0
$ws.Range("H81").Value = 22
# Row 101
$ws.Range("D101").Value = 10
$ws.Range("E101").Value = 391
# Row 119
$ws.Range("B119").Value = 241
$ws.Range("C119").Value = 3
$ws.Range("E119").Value = 164
# Row 144
$ws.Range("B144").Value = 83
$ws.Range("C144").Value = 2
$ws.Range("D144").Value = 48
$ws.Range("E144").Value = 30
# Row 171
$ws.Range("A171").Value = "Maldivas"
$ws.Range("B171").Value = 28
$ws.Range("C171").Value = 3
$ws.Range("D171").Value = 16
$ws.Range("E171").Value = 12
# Row 172
$ws.Range("A172").Value = "Republica del Chad"
$ws.Range("B172").Value = 27
$ws.Range("D172").Value = 5
$ws.Range("E172").Value = 22
# Row 177
$ws.Range("A177").Value = "Timor Oriental"
$ws.Range("F177").Value = 0
# Row 178
$ws.Range("A178").Value = "Belice"
$ws.Range("D178").Value = 0
$ws.Range("E178").Value = 16
$ws.Range("F178").Value = 1
$ws.Range("H178").Value = 2
# Row 179
$ws.Range("A179").Value = "Nueva Caledonia"
$ws.Range("D179").Value = 14
$ws.Range("E179").Value = 4
$ws.Range("H179").Value = 0
